$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the horizontal centering from the "Band ID" cell style while
#        the three original merged blocks still exist (keeps vertical
#        center + wrap text, matches styles.xml change). ---
$ws.Range("C1:D5").HorizontalAlignment = 1
$ws.Range("C7:D10").HorizontalAlignment = 1
$ws.Range("C12:D18").HorizontalAlignment = 1

# --- 2. Unmerge the three "Band ID" blocks (mergeCells removed entirely). ---
$ws.Range("C1:D5").UnMerge()
$ws.Range("C7:D10").UnMerge()
$ws.Range("C12:D18").UnMerge()

# --- 3. Drop the old "#...Band ID: ..." text values from column C (rows
#        1, 7, 12) -- the new data no longer carries band-id labels. ---
$ws.Range("C1").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("C12").ClearContents()

# --- 4. Remove the trailing rows (15-18) that disappear in the new table. ---
$ws.Rows("15:18").Delete()

# --- 5. Clear out the C/D "band id" placeholder cells from row 10 onward --
#        only rows 1-9 keep the styled (empty) C/D cells in the new layout.
$ws.Range("C10:D14").Clear()

# --- 6. Rewrite column A (the instrument list) with the new values,
#        inserting the two rows that used to be blank gaps (6 and 11) so
#        the sheet becomes a contiguous A1:A14 block. ---
$instruments = @(
    " 'Saxophone'",
    " 'Piano'",
    " 'Double bass'",
    " 'Drums'",
    " 'Guitar'",
    " 'Bass'",
    " 'Vocals'",
    " 'Electric guitar'",
    " 'Acoustic guitar'",
    " 'Drum kit'",
    " 'Harmonica'",
    " 'Maraces'",
    " 'Percussion'",
    " 'Electric Harp'"
)

for ($i = 0; $i -lt $instruments.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value2 = $instruments[$i]
}

# --- 7. Give rows 1-9 the (now general-aligned) style on C/D so they match
#        the rest of the styled-but-empty placeholder cells. ---
$ws.Range("C6:D6").Style = $ws.Range("C5:D5").Style

# --- 8. Selection / active cell mirrors the author's final click (A12:A14). ---
$ws.Range("A12:A14").Select()
$excel.ActiveWindow.RangeSelection.Parent.Application.ActiveCell
$ws.Cells.Item(12,1).Activate()

Write-Output "done"
